$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 371 (pushes existing rows 371-420 down to 372-421),
# carrying over the D-column date style (row insert already copies formats
# from the row above, matching Excel's native "Insert Copied Cells" / shift
# behaviour used when a new weekly record is added to the top of this block).
$ws.Rows.Item(371).Insert()

# Populate the newly inserted row 371 with the new weekly record.
$ws.Cells.Item(371, 1).Value = 8
$ws.Cells.Item(371, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(371, 3).Value = 'Coquimbo'
$ws.Cells.Item(371, 4).Value = 45127
$ws.Cells.Item(371, 5).Value = 4
$ws.Cells.Item(371, 6).Value = 100112031
$ws.Cells.Item(371, 7).Value = 'Poroto verde'
$ws.Cells.Item(371, 8).Value = 'Magnum'
$ws.Cells.Item(371, 9).Value = 'Primera'
$ws.Cells.Item(371, 10).Value = 400
$ws.Cells.Item(371, 11).Value = 26000
$ws.Cells.Item(371, 12).Value = 27000
$ws.Cells.Item(371, 13).Value = 26500
$ws.Cells.Item(371, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(371, 15).Value = 'Perú'
$ws.Cells.Item(371, 16).Value = 1060
$ws.Cells.Item(371, 17).Value = 25
$ws.Cells.Item(371, 18).Value = 'Hortaliza'
